$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.385.48'
$ws.Range('E2').Value = '  +4.29%  '
$ws.Range('D3').Value = '1.717.55'
$ws.Range('E3').Value = '  +1.69%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.522'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.74'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.267'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0632'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0893'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '1.964.70'
$ws.Range('E12').Value = '  +1.98%  '
$ws.Range('D13').Value = '1.725.00'
$ws.Range('E13').Value = '  +2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.560'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '28.361.16'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '246.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.11%  '
$ws.Range('D19').Value = '0.0₃0744'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.57'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0512'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').Value = '1.473.28'
$ws.Range('E33').Value = '  -5.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.969'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.99%  '
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.597'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.40'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  -2.25%  '
$ws.Range('D44').Value = '1.869.57'
$ws.Range('E44').Value = '  +1.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.804'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').Value = '0.0₆0109'
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.08'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.103'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.62%  '
